# Animals.xlsx — add 6 new rows (19-24) of animal/procedure records to
# Tabelle1, matching the rows already present (same column layout:
# A=id, B=sex, C=procedure, D=group, E=duration).
#
# xlPasteFormats = -4122 (used to carry over the number/cell formatting
# from the existing rows without clobbering the values we just typed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 19: 85, M, EPI, V, 7 (fully formatted like the preceding rows) ---
$ws.Range("A19").Value = 85
$ws.Range("B19").Value = "M"
$ws.Range("C19").Value = "EPI"
$ws.Range("D19").Value = "V"
$ws.Range("E19").Value = 7
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)

# --- Row 20: 86, M, EPI, V, 7 ---
$ws.Range("A20").Value = 86
$ws.Range("B20").Value = "M"
$ws.Range("C20").Value = "EPI"
$ws.Range("D20").Value = "V"
$ws.Range("E20").Value = 7
$ws.Range("A18:E18").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)

# --- Row 21: 87, M, EPI, V, 7 ---
$ws.Range("A21").Value = 87
$ws.Range("B21").Value = "M"
$ws.Range("C21").Value = "EPI"
$ws.Range("D21").Value = "V"
$ws.Range("E21").Value = 7
$ws.Range("A18:E18").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)

# --- Row 22: 88, M, EPI, V+P, 7 (only A:D carry over the cell style) ---
$ws.Range("A22").Value = 88
$ws.Range("B22").Value = "M"
$ws.Range("C22").Value = "EPI"
$ws.Range("D22").Value = "V+P"
$ws.Range("E22").Value = 7
$ws.Range("A18:D18").Copy()
$ws.Range("A22:D22").PasteSpecial(-4122)

# --- Row 23: 89, M, EPI, F, 7 (only A:C carry over the cell style) ---
$ws.Range("A23").Value = 89
$ws.Range("B23").Value = "M"
$ws.Range("C23").Value = "EPI"
$ws.Range("D23").Value = "F"
$ws.Range("E23").Value = 7
$ws.Range("A18:C18").Copy()
$ws.Range("A23:C23").PasteSpecial(-4122)

# --- Row 24: 90, M, EPI, tachosil only, 7 (only A:D carry over the cell style) ---
$ws.Range("A24").Value = 90
$ws.Range("B24").Value = "M"
$ws.Range("C24").Value = "EPI"
$ws.Range("D24").Value = "tachosil only"
$ws.Range("E24").Value = 7
$ws.Range("A18:D18").Copy()
$ws.Range("A24:D24").PasteSpecial(-4122)

# Match the row heights of the rest of the table (15.75pt).
$ws.Rows.Item(19).RowHeight = 15.75
$ws.Rows.Item(20).RowHeight = 15.75
$ws.Rows.Item(21).RowHeight = 15.75
$ws.Rows.Item(22).RowHeight = 15.75
$ws.Rows.Item(23).RowHeight = 15.75
$ws.Rows.Item(24).RowHeight = 15.75

# Scroll/selection as left by the editor after entering the new data.
$ws.Range("E19:E23").Select()
